$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Functions")

# Mark the POWER function (row 8) and the LN/LG/LOG/LOGN functions (rows 20-23)
# as "IP" (In Progress) instead of "N/A" now that n-param JSON function parsing
# landed for them.
$ws.Range("D8").Value = "IP"
$ws.Range("D20").Value = "IP"
$ws.Range("D21").Value = "IP"
$ws.Range("D22").Value = "IP"
$ws.Range("D23").Value = "IP"

# Widen column B (function names) to fit the longer status text / content.
$ws.Columns.Item(2).ColumnWidth = 25.6667

# Reflect the author's last selection on the sheet.
$ws.Range("H8").Select()
